$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.010394811630249
$ws.Range("B1").Value = 2.429744958877563
$ws.Range("C1").Value = 2.546634912490845
$ws.Range("D1").Value = 3.19565749168396
$ws.Range("E1").Value = 1.230603933334351
